$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.831.54'
$ws.Range("E2").Value = '  +2.78%  '
$ws.Range("D3").Value = '1.868.24'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.94'
$ws.Range("E5").Value = '  +3.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7014'
$ws.Range("E6").Value = '  +2.27%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07787'
$ws.Range("E8").Value = '  +2.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3090'
$ws.Range("E9").Value = '  +2.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.95'
$ws.Range("E10").Value = '  +2.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07848'
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.196'
$ws.Range("E12").Value = '  +3.31%  '
$ws.Range("D13").Value = '1.871.07'
$ws.Range("E13").Value = '  +1.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.88'
$ws.Range("E14").Value = '  +2.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6977'
$ws.Range("E15").Value = '  +3.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.655'
$ws.Range("E16").Value = '  +3.52%  '
$ws.Range("D17").Value = '29.823.86'
$ws.Range("E17").Value = '  +2.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008428'
$ws.Range("E18").Value = '  +2.02%  '
$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").Value = '2.118.05'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '244.42'
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.85'
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.671'
$ws.Range("E23").Value = '  +3.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1518'
$ws.Range("E25").Value = '  +2.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.986'
$ws.Range("E26").Value = '  +3.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.40'
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.46'
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.547'
$ws.Range("E29").Value = '  +1.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.295'
$ws.Range("E30").Value = '  +2.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.245'
$ws.Range("E31").Value = '  +2.83%  '
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05134'
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7899'
$ws.Range("E34").Value = '  +5.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.943'
$ws.Range("E35").Value = '  +7.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.168'
$ws.Range("E36").Value = '  +1.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.710'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '1.332.31'
$ws.Range("E38").Value = '  +10.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01890'
$ws.Range("E39").Value = '  +3.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.748'
$ws.Range("E40").Value = '  +1.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9678'
$ws.Range("E41").Value = '  +5.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.083'
$ws.Range("E42").Value = '  +12.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '107.69'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.837'
$ws.Range("E45").Value = '  +4.96%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.020.11'
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.83'
$ws.Range("E47").Value = '  +3.24%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.797'
$ws.Range("E48").Value = '  +4.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5201'
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000120'
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.050'
$ws.Range("E51").Value = '  +2.43%  '
